# Add support for URL's: new "Extra" column (H) on the "attributes" sheet,
# filled with "noval" placeholders for every existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("attributes")

# New header in H1
$ws.Cells.Item(1, 8).Value = "Extra"

# Fill H2:H17 with the "noval" placeholder
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 8).Value = "noval"
}

# Column G (type_class_id values) gets an explicit best-fit-like width
$ws.Columns.Item(7).ColumnWidth = 41

# Move the active selection as recorded in the saved view state
$ws.Range("D31").Select()
